$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.856.37'
$ws.Range('E2').Value = '  +3.11%  '
$ws.Range('D3').Value = '1.724.33'
$ws.Range('E3').Value = '  +3.02%  '
$ws.Range('E4').Value = '  -0.41%  '
$ws.Range('D5').Value = '''217.35'
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('D7').Value = '''0.996'
$ws.Range('E7').Value = '  -0.48%  '
$ws.Range('D8').Value = '''24.07'
$ws.Range('E8').Value = '  +12.22%  '
$ws.Range('D9').Value = '''0.265'
$ws.Range('E9').Value = '  +4.48%  '
$ws.Range('D10').Value = '''0.0631'
$ws.Range('E10').Value = '  +1.30%  '
$ws.Range('D12').Value = '1.967.62'
$ws.Range('E12').Value = '  +2.99%  '
$ws.Range('D13').Value = '1.721.32'
$ws.Range('E13').Value = '  +2.71%  '
$ws.Range('D14').Value = '''4.24'
$ws.Range('E14').Value = '  +3.14%  '
$ws.Range('D15').Value = '''0.565'
$ws.Range('E15').Value = '  +5.79%  '
$ws.Range('D16').Value = '''68.13'
$ws.Range('E16').Value = '  +2.81%  '
$ws.Range('D17').Value = '27.851.81'
$ws.Range('E17').Value = '  +3.09%  '
$ws.Range('D18').Value = '''241.98'
$ws.Range('E18').Value = '  +2.63%  '
$ws.Range('D19').Value = '''8.09'
$ws.Range('E19').Value = '  -1.32%  '
$ws.Range('E20').Value = '  +1.90%  '
$ws.Range('D21').Value = '''0.995'
$ws.Range('D22').Value = '''4.64'
$ws.Range('E22').Value = '  +3.80%  '
$ws.Range('D23').Value = '''9.73'
$ws.Range('E23').Value = '  +4.97%  '
$ws.Range('D24').Value = '''2.12'
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('D25').Value = '''148.63'
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('E26').Value = '  +3.89%  '
$ws.Range('D27').Value = '''16.64'
$ws.Range('E27').Value = '  +1.13%  '
$ws.Range('E28').Value = '  +1.35%  '
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('E30').Value = '  +1.71%  '
$ws.Range('E31').Value = '  +1.32%  '
$ws.Range('E32').Value = '  +2.23%  '
$ws.Range('D33').Value = '''3.31'
$ws.Range('E33').Value = '  +4.57%  '
$ws.Range('D34').Value = '1.482.39'
$ws.Range('E34').Value = '  -3.92%  '
$ws.Range('E35').Value = '  -1.85%  '
$ws.Range('D36').Value = '''0.967'
$ws.Range('E36').Value = '  +5.99%  '
$ws.Range('D37').Value = '''0.615'
$ws.Range('E37').Value = '  +4.20%  '
$ws.Range('D38').Value = '''2.39'
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('E40').Value = '  +2.77%  '
$ws.Range('D41').Value = '''71.68'
$ws.Range('E41').Value = '  +5.86%  '
$ws.Range('D42').Value = '''5.88'
$ws.Range('E42').Value = '  +6.02%  '
$ws.Range('E43').Value = '  -0.43%  '
$ws.Range('E44').Value = '  +1.93%  '
$ws.Range('D45').Value = '1.871.49'
$ws.Range('E45').Value = '  +2.99%  '
$ws.Range('E46').Value = '  +1.52%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = '''1.68'
$ws.Range('E47').Value = '  +9.43%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '''91.94'
$ws.Range('E48').Value = '  +1.47%  '
$ws.Range('D49').Value = '0.0₆0112'
$ws.Range('E49').Value = '  +3.76%  '
$ws.Range('D50').Value = '''8.35'
$ws.Range('E50').Value = '  +4.08%  '
$ws.Range('E51').Value = '  +2.26%  '
